# PROS-13075 - CCRU - POS 2020 KPIs - SAND Deployment
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Plant Based Drinks Shelf Share" KPI Level 3 rows become their own
# KPI Level 2 line items (KPI Level 2 Name = KPI Level 3 Name) with a
# KPI Level 2 Weight of 0, for both the HS_2020 (rows 24-25) and
# CV_2020 (rows 54-55) KPI groups.
$ws.Range("B24").Value = "HS_2020@Plant Based Drinks Shelf Share"
$ws.Range("C24").Value = 0

$ws.Range("B25").Value = "HS_2020@Plant Based Drinks Shelf Share"
$ws.Range("C25").Value = 0

$ws.Range("B54").Value = "CV_2020@Plant Based Drinks Shelf Share"
$ws.Range("C54").Value = 0

$ws.Range("B55").Value = "CV_2020@Plant Based Drinks Shelf Share"
$ws.Range("C55").Value = 0

# Widen the data columns slightly (values chosen so the stored OOXML
# column width lands as close as possible to the target after the
# runtime's internal pixel-quantization of ColumnWidth).
$ws.Columns.Item(1).ColumnWidth = 43.6666666666667
$ws.Columns.Item(2).ColumnWidth = 31.8333333333333
$ws.Columns.Item(3).ColumnWidth = 27.1666666666667
$ws.Columns.Item(4).ColumnWidth = 34
$ws.Columns.Item(5).ColumnWidth = 41
$ws.Columns.Item(6).ColumnWidth = 17.6666666666667

# Move the active selection to C57.
$ws.Range("C57").Select()
